$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Post" table gains a new attribute column for the actual image bytes:
# a "PostData" attribute of type "blob", inserted right after "World Visible"
# and before the old "JunctionTagID" column. Inserting the column shifts the
# Tag Junction Table / Comment Junction Table blocks one column to the right.
$ws.Range("F1:F3").EntireColumn.Insert()

# The freshly inserted column is blank; give it the same border formatting
# the old "JunctionTagID" column (now shifted into G) had, so the table
# grid lines stay continuous.
$ws.Range("G2:G3").Copy()
$ws.Range("F2:F3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F2").Value = "PostData"
$ws.Range("F3").Value = "blob"

# Resize the new column to fit its contents, like the rest of the table.
$ws.Range("F1:F3").EntireColumn.AutoFit()

$ws.Cells.Item(4, 6).Select()
